$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1837270341207349
$ws.Range("C2").Value = 0.5590551181102362
$ws.Range("J2").Value = 0.02099737532808399
$ws.Range("P2").Value = 0.1679790026246719
$ws.Range("S2").Value = 0.06824146981627296
$ws.Range("B3").Value = 0.01382488479262673
$ws.Range("C3").Value = 0.0184331797235023
$ws.Range("J3").Value = 0.0184331797235023
$ws.Range("P3").Value = 0.7880184331797235
$ws.Range("S3").Value = 0.1612903225806452
$ws.Range("J4").Value = 0.1063829787234043
$ws.Range("P4").Value = 0.574468085106383
$ws.Range("S4").Value = 0.3191489361702128
$ws.Range("B6").Value = 0.06866952789699571
$ws.Range("D6").Value = 0.004291845493562232
$ws.Range("E6").Value = 0.004291845493562232
$ws.Range("F6").Value = 0.04291845493562232
$ws.Range("J6").Value = 0.3090128755364807
$ws.Range("O6").Value = 0.02145922746781116
$ws.Range("Q6").Value = 0.1502145922746781
$ws.Range("R6").Value = 0.07725321888412018
$ws.Range("S6").Value = 0.3218884120171674
$ws.Range("B7").Value = 0.1300813008130081
$ws.Range("D7").Value = 0.01219512195121951
$ws.Range("F7").Value = 0.04878048780487805
$ws.Range("J7").Value = 0.1056910569105691
$ws.Range("O7").Value = 0.04471544715447155
$ws.Range("Q7").Value = 0.1382113821138211
$ws.Range("R7").Value = 0.07317073170731707
$ws.Range("S7").Value = 0.4471544715447154
$ws.Range("B8").Value = 0.1077504725897921
$ws.Range("D8").Value = 0.0113421550094518
$ws.Range("F8").Value = 0.0661625708884688
$ws.Range("J8").Value = 0.109640831758034
$ws.Range("O8").Value = 0.02646502835538752
$ws.Range("Q8").Value = 0.1398865784499055
$ws.Range("R8").Value = 0.1209829867674858
$ws.Range("S8").Value = 0.4177693761814745
$ws.Range("B9").Value = 0.1129032258064516
$ws.Range("D9").Value = 0.02150537634408602
$ws.Range("F9").Value = 0.08064516129032258
$ws.Range("J9").Value = 0.1075268817204301
$ws.Range("O9").Value = 0.03225806451612903
$ws.Range("Q9").Value = 0.1397849462365591
$ws.Range("R9").Value = 0.05376344086021505
$ws.Range("S9").Value = 0.4516129032258064
$ws.Range("B10").Value = 0.125
$ws.Range("D10").Value = 0.02631578947368421
$ws.Range("E10").Value = 0.0007309941520467836
$ws.Range("F10").Value = 0.0577485380116959
$ws.Range("J10").Value = 0.1264619883040936
$ws.Range("O10").Value = 0.01827485380116959
$ws.Range("Q10").Value = 0.1988304093567251
$ws.Range("R10").Value = 0.08406432748538012
$ws.Range("S10").Value = 0.3625730994152047
$ws.Range("G11").Value = 0.1596009975062344
$ws.Range("J11").Value = 0.08478802992518704
$ws.Range("K11").Value = 0.2144638403990025
$ws.Range("L11").Value = 0.5236907730673317
$ws.Range("S11").Value = 0.01745635910224439
$ws.Range("G12").Value = 0.7083333333333334
$ws.Range("J12").Value = 0.2407407407407407
$ws.Range("K12").Value = 0.009259259259259259
$ws.Range("L12").Value = 0.01851851851851852
$ws.Range("S12").Value = 0.02314814814814815
$ws.Range("G13").Value = 0.7291666666666666
$ws.Range("J13").Value = 0.1875
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("H15").Value = 0.1134453781512605
$ws.Range("I15").Value = 0.09243697478991597
$ws.Range("J15").Value = 0.3865546218487395
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("M15").Value = 0.008403361344537815
$ws.Range("O15").Value = 0.03781512605042017
$ws.Range("S15").Value = 0.3067226890756303
$ws.Range("F16").Value = 0.0196078431372549
$ws.Range("H16").Value = 0.1882352941176471
$ws.Range("I16").Value = 0.08235294117647059
$ws.Range("J16").Value = 0.4117647058823529
$ws.Range("K16").Value = 0.1372549019607843
$ws.Range("M16").Value = 0.02352941176470588
$ws.Range("O16").Value = 0.04705882352941176
$ws.Range("S16").Value = 0.09019607843137255
$ws.Range("F17").Value = 0.02546296296296296
$ws.Range("H17").Value = 0.1875
$ws.Range("I17").Value = 0.09490740740740741
$ws.Range("J17").Value = 0.3796296296296297
$ws.Range("K17").Value = 0.1064814814814815
$ws.Range("M17").Value = 0.01851851851851852
$ws.Range("N17").Value = 0.002314814814814815
$ws.Range("O17").Value = 0.05787037037037037
$ws.Range("S17").Value = 0.1273148148148148
$ws.Range("F18").Value = 0.03539823008849557
$ws.Range("H18").Value = 0.2123893805309734
$ws.Range("I18").Value = 0.08849557522123894
$ws.Range("J18").Value = 0.3849557522123894
$ws.Range("K18").Value = 0.0752212389380531
$ws.Range("M18").Value = 0.02212389380530973
$ws.Range("O18").Value = 0.05752212389380531
$ws.Range("S18").Value = 0.1238938053097345
$ws.Range("F19").Value = 0.02407932011331445
$ws.Range("H19").Value = 0.231586402266289
$ws.Range("I19").Value = 0.0594900849858357
$ws.Range("J19").Value = 0.3420679886685553
$ws.Range("K19").Value = 0.1373937677053824
$ws.Range("M19").Value = 0.01912181303116147
$ws.Range("O19").Value = 0.06303116147308782
$ws.Range("S19").Value = 0.1232294617563739
